$d = $word.ActiveDocument

# 1. Header phone number: "+61 0452-614-807" -> "0348-338-392"
$d.Content.Find.Execute("+61 0452-614-807", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "0348-338-392", 2)

# 2. OBJECTIVE paragraph: "Data scientist" -> "Data Scientist" (capitalize S)
$d.Content.Find.Execute("Data scientist in marketing", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Data Scientist in marketing", 2)

# 3. Technical Skills paragraph: add "Data Analytics, " and "Data Ming, "
$d.Content.Find.Execute("Programming, Machine Learning", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Programming, Data Analytics, Machine Learning", 2)

$d.Content.Find.Execute("deep learning), Databases", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "deep learning), Data Ming, Databases", 2)

# 4. Tools paragraph: "MySQL" -> "NoSQL"
$d.Content.Find.Execute("SQL, MySQL, PostgreSQL", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "SQL, NoSQL, PostgreSQL", 2)

# 5. Zalo / Data Scientist dates: "7/2020 - 10/2020" -> "06/2020 - 10/2020"
# (search stays inside the plain italic date run, not the bold job-title run)
$d.Content.Find.Execute("7/2020 – 10/2020", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "06/2020 – 10/2020", 2)

# 6. TMA Solutions job title: "Backend Developer Intern" -> "Data Engineer Intern"
$d.Content.Find.Execute("Backend Developer Intern", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Data Engineer Intern", 2)

# 7. TMA Solutions dates: "7/2017 - 9/2017" -> "07/2017 - 09/2017"
# (search stays inside the plain italic date run, not the bold job-title run)
$d.Content.Find.Execute("7/2017 – 9/2017", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "07/2017 – 09/2017", 2)
